# V1 primera carga basica
# Re-order / update the GPIO-Teclado-Macros table contents and refresh the
# row-height layout + selection to match the latest save of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New button/GPIO/action labels introduced by this revision -------------
# (written first, in the order they were entered, so brand-new text shows
# up in the workbook in the same sequence as the source edit)
$ws.Range("D9").Value  = "D21"
$ws.Range("D3").Value  = "D23"
$ws.Range("D12").Value = "D33"
$ws.Range("E9").Value  = "Ctrl + Alt + Shift + F8"
$ws.Range("E10").Value = "Ctrl + Alt + Shift + F9"
$ws.Range("E11").Value = "Ctrl + Alt + Shift + F10"
$ws.Range("E14").Value = "CTR C"
$ws.Range("E13").Value = "CTR V"
$ws.Range("E15").Value = "CTR Z"
$ws.Range("E16").Value = "Page UP"
$ws.Range("E12").Value = "Ctrl + Alt + Shift + F11"

# --- Remaining cells re-shuffled among the existing button/GPIO values -----
$ws.Range("C3").Value  = "S8"
$ws.Range("C4").Value  = "S7"
$ws.Range("D4").Value  = "D22"
$ws.Range("C6").Value  = "S2"
$ws.Range("D6").Value  = "D16"
$ws.Range("C7").Value  = "S1"
$ws.Range("D7").Value  = "D4"
$ws.Range("C8").Value  = "S13"
$ws.Range("D8").Value  = "D25"
$ws.Range("C9").Value  = "S6"
$ws.Range("C11").Value = "S3"
$ws.Range("D11").Value = "D17"
$ws.Range("C12").Value = "S14"
$ws.Range("C13").Value = "S15"
$ws.Range("D13").Value = "D32"
$ws.Range("C14").Value = "S12"
$ws.Range("D14").Value = "D26"
$ws.Range("C16").Value = "S10"
$ws.Range("D16").Value = "D14"
$ws.Range("C17").Value = "S9"
$ws.Range("D17").Value = "D13"

# --- Row heights (content re-flowed with the new text) ----------------------
$ws.Rows.Item(3).RowHeight  = 28.8
$ws.Rows.Item(4).RowHeight  = 28.8
$ws.Rows.Item(5).RowHeight  = 28.8
$ws.Rows.Item(6).RowHeight  = 28.8
$ws.Rows.Item(7).RowHeight  = 28.8
$ws.Rows.Item(8).RowHeight  = 28.8
$ws.Rows.Item(9).RowHeight  = 28.8
$ws.Rows.Item(10).RowHeight = 28.8
$ws.Rows.Item(11).RowHeight = 28.8
$ws.Rows.Item(12).RowHeight = 28.8

# --- Selection left at the end of the edit session ---------------------------
$ws.Range("E11:E12").Select()
